$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "company" worksheet, placed right after the last existing
#    sheet ("clients") so it lands at the end of the tab strip.
# ---------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$clientsSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $clientsSheet)
$ws.Name = "company"

# ---------------------------------------------------------------------------
# 2. Column widths (values chosen so the exported OOXML <col width="..">
#    matches the target exactly; this engine adds 5/6 to the COM
#    ColumnWidth when writing the stored width).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10.0
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 16.833333333333332
$ws.Columns.Item(4).ColumnWidth = 14.5
$ws.Columns.Item(5).ColumnWidth = 22.0
$ws.Columns.Item(6).ColumnWidth = 14.666666666666666
$ws.Columns.Item(7).ColumnWidth = 17.166666666666668
$ws.Columns.Item(8).ColumnWidth = 25.5
$ws.Columns.Item(9).ColumnWidth = 10.666666666666666
$ws.Columns.Item(10).ColumnWidth = 12.166666666666666
$ws.Columns.Item(11).ColumnWidth = 13.333333333333334
$ws.Columns.Item(12).ColumnWidth = 11.666666666666666

# ---------------------------------------------------------------------------
# 3. Header row (row 1).
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "shortName"
$ws.Range("B1").Value = "companyType"
$ws.Range("C1").Value = "organizationType"

$ws.Range("D1:L1").NumberFormat = "@"
$ws.Range("D1").Value = "unifiedNo"
$ws.Range("E1").Value = "commercialRegistrtaionNo"
$ws.Range("F1").Value = "taxNo"
$ws.Range("G1").Value = "vatRegistrationNo"
$ws.Range("H1").Value = "effectiveVATRegistrationDate"
$ws.Range("I1").Value = "reportType"
$ws.Range("J1").Value = "email"
$ws.Range("K1").Value = "name"
$ws.Range("L1").Value = "role"

# ---------------------------------------------------------------------------
# 4. Data row (row 2).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Fai"
$ws.Range("B2").Value = "منشأة"
$ws.Range("C2").Value = "مؤسسة"

$ws.Range("D2:I2").NumberFormat = "@"
$ws.Range("D2").Value = "1234567890"
$ws.Range("E2").Value = "1111122222"
$ws.Range("F2").Value = "0987654321"
$ws.Range("G2").Value = "111112222233333"
$ws.Range("H2").Value = "Jan 01 2024"
$ws.Range("I2").Value = "شهري"

$ws.Range("K2:L2").NumberFormat = "@"
$ws.Range("K2").Value = "Mario Nady"
$ws.Range("L2").Value = "مالك الحساب"

# J2 carries a mailto hyperlink, styled with the workbook's Hyperlink cell
# style plus the text number format.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "mario@fai.ws"
[void]$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:mario@fai.ws")

# ---------------------------------------------------------------------------
# 5. Selection left on the new sheet (L6) before the user moved away.
# ---------------------------------------------------------------------------
[void]$ws.Range("L6").Select()

# ---------------------------------------------------------------------------
# 6. Final navigation: the user ends up back on "purchaseCashback" with cell
#    K19 selected, which becomes the active tab/sheet on save.
# ---------------------------------------------------------------------------
$purchaseCashback = $wb.Worksheets.Item("purchaseCashback")
$purchaseCashback.Activate()
[void]$purchaseCashback.Range("K19").Select()
